$d = $word.ActiveDocument

# Helper: locate the whole paragraph (including its end-of-paragraph mark)
# that contains the given search text, using Find to get the starting
# offset and then matching it against the document's Paragraphs collection.
# This is more robust than hard-coded paragraph indices.
function Get-ParaRangeByText($searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $searchText"
    }
    $targetStart = $rng.Start
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $pr = $p.Range
        if ($pr.Start -eq $targetStart) {
            return $pr
        }
    }
    throw "Paragraph containing text not located: $searchText"
}

function Set-ParagraphXml($paraRange, $styleVal, $innerRunsXml) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p><w:pPr><w:pStyle w:val="' + $styleVal + '"/></w:pPr>' + $innerRunsXml + '</w:p>' +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $null = $paraRange.InsertXML($xml)
}

# --- Edit 1: "Цель работы" paragraph --------------------------------------
# Split the old single sentence into two sentences separated by a manual
# line break, keeping the trailing " [1]" runs untouched.
$r1 = Get-ParaRangeByText("Научиться оформлять отчёты с помощью легковесного языка разметки Markdown.")
$runs1 = '<w:r><w:t xml:space="preserve">Изучить идеологию и применение средств контроля версий.</w:t></w:r>' +
    '<w:r><w:br/></w:r>' +
    '<w:r><w:t xml:space="preserve">Освоить умения по работе с git.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">[1]</w:t></w:r>'
Set-ParagraphXml $r1 "FirstParagraph" $runs1

# --- Edit 2: "Задание" paragraph -------------------------------------------
# Replace the whole paragraph (five runs) with a single new run.
$r2 = Get-ParaRangeByText("Сделайте отчёт по предыдущей лабораторной работе в формате Markdown.")
$runs2 = '<w:r><w:t xml:space="preserve">Научиться применять команды git, работать с github.</w:t></w:r>'
Set-ParagraphXml $r2 "FirstParagraph" $runs2

# --- Edit 3: "Выводы" paragraph --------------------------------------------
# Split the old single sentence into two sentences separated by a space run.
$r3 = Get-ParaRangeByText("В ходе работы мы научились оформлять отчёты с помощью легковесного языка разметки Markdown.")
$runs3 = '<w:r><w:t xml:space="preserve">В ходе работы мы изучили идеологию и применение средств контроля версий.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">Освоили умения по работе с git.</w:t></w:r>'
Set-ParagraphXml $r3 "FirstParagraph" $runs3

Write-Host "edits applied"
